# Apply Covid DailyStats update: "po 12. 07. 2021"
# - Revises AgTests (F) / AgPosit (G) figures for many existing rows (306-488)
# - Appends 5 new daily rows (489-493) covering 2021-07-06 .. 2021-07-10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised F (AgTests) / G (AgPosit) values for existing rows ---
$ws.Range("F306").Value = 76845
$ws.Range("G306").Value = 7682
$ws.Range("F307").Value = 75268
$ws.Range("G307").Value = 6312
$ws.Range("F309").Value = 77701
$ws.Range("G309").Value = 5512
$ws.Range("F310").Value = 79517
$ws.Range("F313").Value = 76674
$ws.Range("G313").Value = 3462
$ws.Range("F314").Value = 65281
$ws.Range("G314").Value = 3135
$ws.Range("F315").Value = 56819
$ws.Range("G315").Value = 2627
$ws.Range("F316").Value = 50921
$ws.Range("G316").Value = 2302
$ws.Range("F317").Value = 63901
$ws.Range("G317").Value = 2176
$ws.Range("F320").Value = 74170
$ws.Range("G320").Value = 3377
$ws.Range("F321").Value = 90953
$ws.Range("F322").Value = 110209
$ws.Range("G322").Value = 2334
$ws.Range("F323").Value = 217415
$ws.Range("G323").Value = 3102
$ws.Range("F324").Value = 248649
$ws.Range("G324").Value = 2850
$ws.Range("F325").Value = 769915
$ws.Range("G325").Value = 6467
$ws.Range("F327").Value = 224424
$ws.Range("G327").Value = 2720
$ws.Range("F328").Value = 180648
$ws.Range("F329").Value = 73598
$ws.Range("G329").Value = 1729
$ws.Range("F330").Value = 71520
$ws.Range("G330").Value = 2077
$ws.Range("F331").Value = 154320
$ws.Range("F332").Value = 485019
$ws.Range("G332").Value = 4801
$ws.Range("F334").Value = 192895
$ws.Range("G334").Value = 3500
$ws.Range("F335").Value = 150541
$ws.Range("F336").Value = 81977
$ws.Range("F338").Value = 221312
$ws.Range("F339").Value = 661016
$ws.Range("G339").Value = 5488
$ws.Range("F341").Value = 283422
$ws.Range("F343").Value = 133821
$ws.Range("G343").Value = 2983
$ws.Range("F344").Value = 136254
$ws.Range("F346").Value = 674784
$ws.Range("F349").Value = 159489
$ws.Range("G349").Value = 2755
$ws.Range("F351").Value = 150278
$ws.Range("G351").Value = 2799
$ws.Range("F352").Value = 306714
$ws.Range("F356").Value = 160777
$ws.Range("F357").Value = 138302
$ws.Range("G357").Value = 3010
$ws.Range("F359").Value = 320680
$ws.Range("G359").Value = 3331
$ws.Range("F362").Value = 228968
$ws.Range("F363").Value = 189510
$ws.Range("F366").Value = 338788
$ws.Range("F369").Value = 235139
$ws.Range("F370").Value = 180668
$ws.Range("G370").Value = 2043
$ws.Range("F371").Value = 160283
$ws.Range("G371").Value = 1969
$ws.Range("F372").Value = 178681
$ws.Range("G372").Value = 1857
$ws.Range("F373").Value = 350518
$ws.Range("G373").Value = 2384
$ws.Range("F374").Value = 773153
$ws.Range("G374").Value = 3420
$ws.Range("F375").Value = 349804
$ws.Range("G375").Value = 1840
$ws.Range("F376").Value = 221918
$ws.Range("F377").Value = 176835
$ws.Range("G377").Value = 1826
$ws.Range("F378").Value = 157496
$ws.Range("G378").Value = 1550
$ws.Range("F379").Value = 180800
$ws.Range("G379").Value = 1614
$ws.Range("F380").Value = 345146
$ws.Range("G380").Value = 2027
$ws.Range("F381").Value = 747266
$ws.Range("F382").Value = 356792
$ws.Range("F383").Value = 221304
$ws.Range("F384").Value = 172327
$ws.Range("G384").Value = 1516
$ws.Range("F385").Value = 151050
$ws.Range("G385").Value = 1407
$ws.Range("F386").Value = 183198
$ws.Range("G386").Value = 1363
$ws.Range("F387").Value = 351621
$ws.Range("G387").Value = 1664
$ws.Range("F388").Value = 729806
$ws.Range("G388").Value = 2202
$ws.Range("F391").Value = 177852
$ws.Range("G391").Value = 1209
$ws.Range("F392").Value = 221841
$ws.Range("F393").Value = 308467
$ws.Range("F395").Value = 752721
$ws.Range("F398").Value = 299174
$ws.Range("G398").Value = 1470
$ws.Range("F399").Value = 200553
$ws.Range("F400").Value = 149219
$ws.Range("G400").Value = 765
$ws.Range("F401").Value = 272700
$ws.Range("F402").Value = 722458
$ws.Range("F404").Value = 224127
$ws.Range("F405").Value = 174734
$ws.Range("F406").Value = 171487
$ws.Range("F407").Value = 158329
$ws.Range("F408").Value = 304887
$ws.Range("F409").Value = 708940
$ws.Range("F410").Value = 364822
$ws.Range("F412").Value = 176625
$ws.Range("F413").Value = 149789
$ws.Range("F414").Value = 149245
$ws.Range("F415").Value = 308205
$ws.Range("G415").Value = 695
$ws.Range("F416").Value = 672005
$ws.Range("F417").Value = 343404
$ws.Range("F418").Value = 202392
$ws.Range("F419").Value = 149703
$ws.Range("F420").Value = 138984
$ws.Range("G420").Value = 501
$ws.Range("F421").Value = 153255
$ws.Range("G421").Value = 534
$ws.Range("F422").Value = 298530
$ws.Range("F424").Value = 265930
$ws.Range("F426").Value = 106968
$ws.Range("F433").Value = 86222
$ws.Range("G433").Value = 265
$ws.Range("F434").Value = 79053
$ws.Range("F435").Value = 83395
$ws.Range("F436").Value = 139576
$ws.Range("F437").Value = 162258
$ws.Range("F440").Value = 73015
$ws.Range("F442").Value = 67425
$ws.Range("F444").Value = 100201
$ws.Range("G444").Value = 177
$ws.Range("F447").Value = 64796
$ws.Range("F449").Value = 59853
$ws.Range("F451").Value = 83029
$ws.Range("G451").Value = 113
$ws.Range("F453").Value = 67302
$ws.Range("F454").Value = 50906
$ws.Range("G455").Value = 118
$ws.Range("F456").Value = 48005
$ws.Range("F457").Value = 75615
$ws.Range("G457").Value = 126
$ws.Range("F458").Value = 67816
$ws.Range("F460").Value = 55770
$ws.Range("F461").Value = 43628
$ws.Range("F465").Value = 58418
$ws.Range("F466").Value = 49563
$ws.Range("F468").Value = 40571
$ws.Range("F469").Value = 39129
$ws.Range("F470").Value = 41452
$ws.Range("F471").Value = 62581
$ws.Range("F472").Value = 47975
$ws.Range("F473").Value = 38854
$ws.Range("F474").Value = 43851
$ws.Range("G474").Value = 58
$ws.Range("F475").Value = 34171
$ws.Range("F476").Value = 35180
$ws.Range("F477").Value = 37135
$ws.Range("F478").Value = 51077
$ws.Range("F479").Value = 39882
$ws.Range("F480").Value = 32620
$ws.Range("F481").Value = 43048
$ws.Range("F482").Value = 34539
$ws.Range("F483").Value = 63670
$ws.Range("F484").Value = 8058
$ws.Range("F485").Value = 13502
$ws.Range("F486").Value = 8448
$ws.Range("G486").Value = 7
$ws.Range("F487").Value = 6698
$ws.Range("G487").Value = 9
$ws.Range("F488").Value = 6146
$ws.Range("G488").Value = 8

# --- New rows 489-493 (2021-07-06 .. 2021-07-10) ---
$ws.Range("A489").Value = 44383
$ws.Range("B489").Value = 391780
$ws.Range("C489").Value = 5637
$ws.Range("D489").Value = 45
$ws.Range("E489").Value = 12516
$ws.Range("F489").Value = 11542
$ws.Range("G489").Value = 10
$ws.Range("A490").Value = 44384
$ws.Range("B490").Value = 391813
$ws.Range("C490").Value = 5245
$ws.Range("D490").Value = 33
$ws.Range("E490").Value = 12516
$ws.Range("F490").Value = 9711
$ws.Range("G490").Value = 15
$ws.Range("A491").Value = 44385
$ws.Range("B491").Value = 391852
$ws.Range("C491").Value = 5671
$ws.Range("D491").Value = 39
$ws.Range("E491").Value = 12517
$ws.Range("F491").Value = 8847
$ws.Range("G491").Value = 10
$ws.Range("A492").Value = 44386
$ws.Range("B492").Value = 391893
$ws.Range("C492").Value = 8903
$ws.Range("D492").Value = 41
$ws.Range("E492").Value = 12519
$ws.Range("F492").Value = 11820
$ws.Range("G492").Value = 11
$ws.Range("A493").Value = 44387
$ws.Range("B493").Value = 391918
$ws.Range("C493").Value = 4699
$ws.Range("D493").Value = 25
$ws.Range("E493").Value = 12519
$ws.Range("F493").Value = 6440
$ws.Range("G493").Value = 2
